$d = $word.ActiveDocument

$find = "Datas das campanhas de constelação de botas 2022"
$replace = "Datas das campanhas de 2022 que usam constelação de botas"

$d.Content.Find.Execute($find, $true, $false, $false, $false, $false,
                         $true, 1, $false, $replace, 2)
